$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at position 163 (shifts rows 163:175 down to 164:176)
$ws.Rows.Item(163).Insert()

# Copy formatting of the row that is now just below (row 164, old row 163)
# down onto the freshly inserted (blank) row 163 so styles (e.g. the date
# number format on column D) carry over correctly. Restrict to the used
# columns (A:R) so we don't drag formatting across the whole 16384-column row.
$ws.Range("A164:R164").Copy()
$ws.Range("A163:R163").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row with its data
$ws.Cells.Item(163, 1).Value = 5
$ws.Cells.Item(163, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(163, 3).Value = "Maule"
$ws.Cells.Item(163, 4).Value = 44568
$ws.Cells.Item(163, 5).Value = 7
$ws.Cells.Item(163, 6).Value = 100112045
$ws.Cells.Item(163, 7).Value = "Zapallo"
$ws.Cells.Item(163, 8).Value = "Camote"
$ws.Cells.Item(163, 9).Value = "1a nueva(o)"
$ws.Cells.Item(163, 10).Value = 800
$ws.Cells.Item(163, 11).Value = 300
$ws.Cells.Item(163, 12).Value = 300
$ws.Cells.Item(163, 13).Value = 300
$ws.Cells.Item(163, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(163, 15).Value = "Región del Maule"
$ws.Cells.Item(163, 16).Value = 300
$ws.Cells.Item(163, 17).Value = 1
$ws.Cells.Item(163, 18).Value = "Hortaliza"
